$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2242648.2
$ws.Range("J17").Value = 2319947.8
$ws.Range("L17").Value = 6959843.399999999
$ws.Range("N17").Value = -6960179.399999999
$ws.Range("H38").Value = 4949.758
$ws.Range("J38").Value = 6540
$ws.Range("L38").Value = 19620
$ws.Range("N38").Value = -20364
$ws.Range("H41").Value = 291.5
$ws.Range("I41").Value = 250
$ws.Range("J41").Value = 333
$ws.Range("K41").Value = 250
$ws.Range("L41").Value = 333
$ws.Range("M41").Value = 190
$ws.Range("N41").Value = -1213
$ws.Range("H42").Value = 242
$ws.Range("I42").Value = 99
$ws.Range("K42").Value = 297
$ws.Range("M42").Value = -67
$ws.Range("H43").Value = 5073.25
$ws.Range("I43").Value = 4597.4
$ws.Range("J43").Value = 5866.3335
$ws.Range("K43").Value = 4597.4
$ws.Range("L43").Value = 5866.3335
$ws.Range("M43").Value = -4528.4
$ws.Range("N43").Value = -6004.3335
$ws.Range("H62").Value = 2618.4443
$ws.Range("I62").Value = 1946.375
$ws.Range("K62").Value = 1946.375
$ws.Range("M62").Value = -1322.375
$ws.Range("H65").Value = 2618.4443
$ws.Range("I65").Value = 1946.375
$ws.Range("K65").Value = 9731.875
$ws.Range("M65").Value = -6611.875
$ws.Range("H98").Value = 751.1875
$ws.Range("I98").Value = 751.1875
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 751.1875
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 746.8125
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 751.1875
$ws.Range("I122").Value = 751.1875
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2253.5625
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 196.4375
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H125").Value = 4534.6
$ws.Range("I125").Value = 7186.3335
$ws.Range("J125").Value = 3398.1428
$ws.Range("K125").Value = 64677.0015
$ws.Range("L125").Value = 30583.2852
$ws.Range("M125").Value = -62217.0015
$ws.Range("N125").Value = -35503.2852
$ws.Range("H131").Value = 3798.8333
$ws.Range("I131").Value = 3158.8
$ws.Range("K131").Value = 9476.400000000001
$ws.Range("M131").Value = -4436.400000000001
$ws.Range("H137").Value = 13083.611
$ws.Range("I137").Value = 14639.23
$ws.Range("J137").Value = 9039
$ws.Range("K137").Value = 43917.69
$ws.Range("L137").Value = 27117
$ws.Range("M137").Value = -41367.69
$ws.Range("N137").Value = -32217
$ws.Range("H138").Value = 18184444
$ws.Range("I138").Value = 1072.25
$ws.Range("J138").Value = 32261892
$ws.Range("K138").Value = 3216.75
$ws.Range("L138").Value = 96785676
$ws.Range("M138").Value = 1923.25
$ws.Range("N138").Value = -96795956

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1761.8823
$ws.Range("I97").Value = 1600.2
$ws.Range("K97").Value = 1600.2
$ws.Range("M97").Value = -1104.2
$ws.Range("H132").Value = 34028.543
$ws.Range("I132").Value = 3815.375
$ws.Range("K132").Value = 11446.125
$ws.Range("M132").Value = -8916.125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2966.9333
$ws.Range("J31").Value = 5403.8887
$ws.Range("L31").Value = 5403.8887
$ws.Range("N31").Value = -5993.8887
$ws.Range("H34").Value = 2966.9333
$ws.Range("J34").Value = 5403.8887
$ws.Range("L34").Value = 5403.8887
$ws.Range("N34").Value = -5807.8887
$ws.Range("H132").Value = 3545.5532
$ws.Range("I132").Value = 3351.125
$ws.Range("J132").Value = 4656.5713
$ws.Range("K132").Value = 10053.375
$ws.Range("L132").Value = 13969.7139
$ws.Range("M132").Value = -7523.375
$ws.Range("N132").Value = -19029.7139
$ws.Range("H134").Value = 2600.818
$ws.Range("I134").Value = 2601.6875
$ws.Range("J134").Value = 2598.5
$ws.Range("K134").Value = 7805.0625
$ws.Range("L134").Value = 7795.5
$ws.Range("M134").Value = -5270.0625
$ws.Range("N134").Value = -12865.5
$ws.Range("H141").Value = 90724.44500000001
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 90724.44500000001
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 90724.44500000001
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -101084.445

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").ClearContents()
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
$ws.Range("H131").Value = 1755.3077
$ws.Range("J131").Value = 1776.2727
$ws.Range("L131").Value = 5328.8181
$ws.Range("N131").Value = -15408.8181

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 632.1818
$ws.Range("I107").Value = 544.375
$ws.Range("K107").Value = 544.375
$ws.Range("M107").Value = 1375.625
$ws.Range("H122").Value = 2368.238
$ws.Range("I122").Value = 1364.4166
$ws.Range("J122").Value = 3706.6667
$ws.Range("K122").Value = 4093.2498
$ws.Range("L122").Value = 11120.0001
$ws.Range("M122").Value = -1643.2498
$ws.Range("N122").Value = -16020.0001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4550.3
$ws.Range("I122").Value = 3188.5
$ws.Range("J122").Value = 9997.5
$ws.Range("K122").Value = 9565.5
$ws.Range("L122").Value = 29992.5
$ws.Range("M122").Value = -7115.5
$ws.Range("N122").Value = -34892.5
$ws.Range("H130").Value = 94977.5
$ws.Range("J130").Value = 94977.5
$ws.Range("L130").Value = 94977.5
$ws.Range("N130").Value = -105017.5
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 2128.9722
$ws.Range("I132").Value = 1988.2759
$ws.Range("K132").Value = 5964.8277
$ws.Range("M132").Value = -3434.8277

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3335.75
$ws.Range("I136").Value = 2353.372
$ws.Range("J136").Value = 5820.5884
$ws.Range("K136").Value = 7060.116
$ws.Range("L136").Value = 17461.7652
$ws.Range("M136").Value = -4510.116
$ws.Range("N136").Value = -22561.7652
